# Auto-generated edit script: update Kujata_Profits market-price derived values
# across all 8 profession sheets, per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 881.17645
$ws.Range("I80").Value = 1182
$ws.Range("J80").Value = 717.0909
$ws.Range("K80").Value = 3546
$ws.Range("L80").Value = 2151.2727
$ws.Range("M80").Value = -2548
$ws.Range("N80").Value = -4147.2727
$ws.Range("H83").Value = 881.17645
$ws.Range("I83").Value = 1182
$ws.Range("J83").Value = 717.0909
$ws.Range("K83").Value = 10638
$ws.Range("L83").Value = 6453.8181
$ws.Range("M83").Value = -5646
$ws.Range("N83").Value = -16437.8181
$ws.Range("H107").Value = 1801.4828
$ws.Range("I107").Value = 1185.84
$ws.Range("K107").Value = 1185.84
$ws.Range("M107").Value = 734.1600000000001
$ws.Range("H111").Value = 980
$ws.Range("I111").Value = 980
$ws.Range("K111").Value = 2940
$ws.Range("M111").Value = 127
$ws.Range("H123").Value = 38749.75
$ws.Range("J123").Value = 38749.75
$ws.Range("L123").Value = 38749.75
$ws.Range("N123").Value = -48549.75
$ws.Range("H129").Value = 815.1875
$ws.Range("I129").Value = 577
$ws.Range("J129").Value = 923.4545000000001
$ws.Range("K129").Value = 1731
$ws.Range("L129").Value = 2770.3635
$ws.Range("M129").Value = 3269
$ws.Range("N129").Value = -12770.3635
$ws.Range("H132").Value = 15887999
$ws.Range("I132").Value = 19618122
$ws.Range("K132").Value = 58854366
$ws.Range("M132").Value = -58851836
$ws.Range("H134").Value = 35976.9
$ws.Range("J134").Value = 35976.9
$ws.Range("L134").Value = 35976.9
$ws.Range("N134").Value = -46116.9
$ws.Range("H135").Value = 43479330
$ws.Range("I135").Value = 796.1429000000001
$ws.Range("K135").Value = 7165.2861
$ws.Range("M135").Value = -4630.2861
$ws.Range("H138").Value = 2084.6555
$ws.Range("I138").Value = 1374.2
$ws.Range("J138").Value = 2287.6428
$ws.Range("K138").Value = 4122.6
$ws.Range("L138").Value = 6862.928400000001
$ws.Range("M138").Value = 1017.4
$ws.Range("N138").Value = -17142.9284

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8094.57
$ws.Range("I32").Value = 5742.9883
$ws.Range("J32").Value = 30042.666
$ws.Range("K32").Value = 5742.9883
$ws.Range("L32").Value = 30042.666
$ws.Range("M32").Value = -5455.9883
$ws.Range("N32").Value = -30616.666
$ws.Range("H43").Value = 12000
$ws.Range("J43").Value = 12000
$ws.Range("L43").Value = 12000
$ws.Range("N43").Value = -12626
$ws.Range("H61").Value = 1417.3939
$ws.Range("I61").Value = 1198.8889
$ws.Range("J61").Value = 2400.6667
$ws.Range("K61").Value = 1198.8889
$ws.Range("L61").Value = 2400.6667
$ws.Range("M61").Value = -986.8888999999999
$ws.Range("N61").Value = -2824.6667
$ws.Range("H63").Value = 2769.647
$ws.Range("I63").Value = 2368.9333
$ws.Range("J63").Value = 5775
$ws.Range("K63").Value = 2368.9333
$ws.Range("L63").Value = 5775
$ws.Range("M63").Value = -1682.9333
$ws.Range("N63").Value = -7147
$ws.Range("H66").Value = 2769.647
$ws.Range("I66").Value = 2368.9333
$ws.Range("J66").Value = 5775
$ws.Range("K66").Value = 11844.6665
$ws.Range("L66").Value = 28875
$ws.Range("M66").Value = -8412.666500000001
$ws.Range("N66").Value = -35739
$ws.Range("H136").Value = 1417.3939
$ws.Range("I136").Value = 1198.8889
$ws.Range("J136").Value = 2400.6667
$ws.Range("K136").Value = 3596.6667
$ws.Range("L136").Value = 7202.000100000001
$ws.Range("M136").Value = -1046.6667
$ws.Range("N136").Value = -12302.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H69").Value = 9833.333000000001
$ws.Range("J69").Value = 9833.333000000001
$ws.Range("L69").Value = 9833.333000000001
$ws.Range("N69").Value = -11455.333
$ws.Range("H72").Value = 9833.333000000001
$ws.Range("J72").Value = 9833.333000000001
$ws.Range("L72").Value = 29499.999
$ws.Range("N72").Value = -37611.999
$ws.Range("H82").Value = 15484.333
$ws.Range("I82").Value = 2025.125
$ws.Range("K82").Value = 2025.125
$ws.Range("M82").Value = -1642.125
$ws.Range("H85").Value = 15484.333
$ws.Range("I85").Value = 2025.125
$ws.Range("K85").Value = 2025.125
$ws.Range("M85").Value = -699.125
$ws.Range("H99").Value = 47620184
$ws.Range("I99").Value = 76924120
$ws.Range("J99").Value = 1294.375
$ws.Range("K99").Value = 76924120
$ws.Range("L99").Value = 1294.375
$ws.Range("M99").Value = -76922622
$ws.Range("N99").Value = -4290.375
$ws.Range("H107").Value = 859.7742
$ws.Range("I107").Value = 643.8461
$ws.Range("J107").Value = 1982.6
$ws.Range("K107").Value = 643.8461
$ws.Range("L107").Value = 1982.6
$ws.Range("M107").Value = 1276.1539
$ws.Range("N107").Value = -5822.6
$ws.Range("H126").Value = 48550
$ws.Range("J126").Value = 48550
$ws.Range("L126").Value = 48550
$ws.Range("N126").Value = -58430

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1868.125
$ws.Range("I31").Value = 1828.3158
$ws.Range("K31").Value = 1828.3158
$ws.Range("M31").Value = -1533.3158
$ws.Range("H34").Value = 1868.125
$ws.Range("I34").Value = 1828.3158
$ws.Range("K34").Value = 1828.3158
$ws.Range("M34").Value = -1626.3158
$ws.Range("H74").Value = 29250
$ws.Range("I74").Value = 17000
$ws.Range("K74").Value = 17000
$ws.Range("M74").Value = -16126
$ws.Range("H77").Value = 29250
$ws.Range("I77").Value = 17000
$ws.Range("K77").Value = 51000
$ws.Range("M77").Value = -46632
$ws.Range("H94").Value = 1999.375
$ws.Range("I94").Value = 1578.6666
$ws.Range("K94").Value = 1578.6666
$ws.Range("M94").Value = -1127.6666
$ws.Range("H108").Value = 33624.8
$ws.Range("J108").Value = 33624.8
$ws.Range("L108").Value = 33624.8
$ws.Range("N108").Value = -41304.8
$ws.Range("H141").Value = 618654
$ws.Range("J141").Value = 685171.1
$ws.Range("L141").Value = 685171.1
$ws.Range("N141").Value = -695531.1

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16418870
$ws.Range("I131").Value = 71429080
$ws.Range("J131").Value = 32848.766
$ws.Range("K131").Value = 214287240
$ws.Range("L131").Value = 98546.29800000001
$ws.Range("M131").Value = -214282200
$ws.Range("N131").Value = -108626.298
$ws.Range("H139").Value = 2968.611
$ws.Range("I139").Value = 3341.4614
$ws.Range("J139").Value = 1999.2
$ws.Range("K139").Value = 10024.3842
$ws.Range("L139").Value = 5997.6
$ws.Range("M139").Value = -4884.3842
$ws.Range("N139").Value = -16277.6

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 327.84616
$ws.Range("I2").Value = 233.66667
$ws.Range("J2").Value = 456.27274
$ws.Range("K2").Value = 233.66667
$ws.Range("L2").Value = 456.27274
$ws.Range("M2").Value = -120.66667
$ws.Range("N2").Value = -682.27274
$ws.Range("H70").Value = 21432608
$ws.Range("I70").Value = 15629242
$ws.Range("J70").Value = 40003380
$ws.Range("K70").Value = 15629242
$ws.Range("L70").Value = 40003380
$ws.Range("M70").Value = -15628972
$ws.Range("N70").Value = -40003920
$ws.Range("H73").Value = 21432608
$ws.Range("I73").Value = 15629242
$ws.Range("J73").Value = 40003380
$ws.Range("K73").Value = 15629242
$ws.Range("L73").Value = 40003380
$ws.Range("M73").Value = -15628306
$ws.Range("N73").Value = -40005252
$ws.Range("H97").Value = 1060.8
$ws.Range("I97").Value = 888.5
$ws.Range("J97").Value = 1750
$ws.Range("K97").Value = 888.5
$ws.Range("L97").Value = 1750
$ws.Range("M97").Value = -392.5
$ws.Range("N97").Value = -2742
$ws.Range("H107").Value = 344.2143
$ws.Range("I107").Value = 363.33334
$ws.Range("J107").Value = 309.8
$ws.Range("K107").Value = 363.33334
$ws.Range("L107").Value = 309.8
$ws.Range("M107").Value = 1556.66666
$ws.Range("N107").Value = -4149.8
$ws.Range("H113").Value = 1114.4615
$ws.Range("I113").Value = 1329.6
$ws.Range("K113").Value = 1329.6
$ws.Range("M113").Value = 840.4000000000001
$ws.Range("H126").Value = 2181.2856
$ws.Range("I126").Value = 1804.3334
$ws.Range("J126").Value = 2859.8
$ws.Range("K126").Value = 5413.0002
$ws.Range("L126").Value = 8579.400000000001
$ws.Range("M126").Value = -2943.0002
$ws.Range("N126").Value = -13519.4
$ws.Range("H128").Value = 39268
$ws.Range("J128").Value = 40780
$ws.Range("L128").Value = 40780
$ws.Range("N128").Value = -50740

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H132").Value = 3386.25
$ws.Range("I132").Value = 3484.75
$ws.Range("J132").Value = 3287.75
$ws.Range("K132").Value = 10454.25
$ws.Range("L132").Value = 9863.25
$ws.Range("M132").Value = -7924.25
$ws.Range("N132").Value = -14923.25
$ws.Range("H136").Value = 1466.45
$ws.Range("I136").Value = 1095.25
$ws.Range("K136").Value = 3285.75
$ws.Range("M136").Value = -735.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H122").Value = 10002115
$ws.Range("I122").Value = 11365722
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 34097166
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -34094716
$ws.Range("N122").Value = -11899.9999

